# Applies the "Trade #393 closed / 4 new trades opened" update described by the diff.
$wb = $excel.ActiveWorkbook

# Helper: write a literal text value into a cell without Excel's "looks like a
# date" auto-conversion kicking in. Only the Date column ("2026-02-18") needs
# this; HH:MM:SS time strings are stored as plain text automatically.
function Set-DateTextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.36   # Current Capital
$summary.Range("B4").Value = 0.47      # Total P&L $
$summary.Range("B5").Value = 0.02      # Total P&L %
$summary.Range("B6").Value = 393       # Total Trades
$summary.Range("B8").Value = 143       # Losing Trades
$summary.Range("B9").Value = 38.93     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - "momentum" row (row 11)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C11").Value = 99.31999999999999
$status.Range("D11").Value = 71
$status.Range("E11").Value = -0.68
$status.Range("F11").Value = -0.68
$status.Range("G11").Value = 32.39

# ---------------------------------------------------------------------------
# All Trades sheet - close trade #393 (row 394) + append trades #422-#425
# Columns: A Trade#  B Date  C Time  D Strategy  E Side  F Entry Price
#          G Exit Price  H Status  I P&L%  J P&L$  K Capital After
#          L Exit Reason  M Duration (min)  N Entry Slippage  O Exit Slippage
#          P Confidence  Q Entry Reason
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #393 (momentum, row 394) transitions from OPEN to CLOSED.
$allTrades.Cells.Item(394, 7).Value = 0.74                    # Exit Price
$allTrades.Cells.Item(394, 8).Value = "CLOSED"                # Status
$allTrades.Cells.Item(394, 9).Value = -7.5                    # P&L %
$allTrades.Cells.Item(394, 10).Value = -0.06                  # P&L $
$allTrades.Cells.Item(394, 11).Value = 99.31999999999999      # Capital After
$allTrades.Cells.Item(394, 12).Value = "early_exit"           # Exit Reason
$allTrades.Cells.Item(394, 13).Value = 0.36                   # Duration (min)

# New trade rows 423-426 (trade #422-#425), all freshly opened.
$allTrades.Cells.Item(423, 1).Value = 422
Set-DateTextCell $allTrades 423 2 "2026-02-18"
$allTrades.Cells.Item(423, 3).Value = "02:03:39"
$allTrades.Cells.Item(423, 4).Value = "momentum"
$allTrades.Cells.Item(423, 5).Value = "DOWN"
$allTrades.Cells.Item(423, 6).Value = 0.8
$allTrades.Cells.Item(423, 8).Value = "OPEN"
$allTrades.Cells.Item(423, 9).Value = 0
$allTrades.Cells.Item(423, 10).Value = 0
$allTrades.Cells.Item(423, 11).Value = 99.37699700270591
$allTrades.Cells.Item(423, 13).Value = 0
$allTrades.Cells.Item(423, 14).Value = 0
$allTrades.Cells.Item(423, 15).Value = 0
$allTrades.Cells.Item(423, 16).Value = 0.9
$allTrades.Cells.Item(423, 17).Value = "Downward momentum: -20.588% over 10 samples"

$allTrades.Cells.Item(424, 1).Value = 423
Set-DateTextCell $allTrades 424 2 "2026-02-18"
$allTrades.Cells.Item(424, 3).Value = "02:03:40"
$allTrades.Cells.Item(424, 4).Value = "HighProbConvergence"
$allTrades.Cells.Item(424, 5).Value = "UP"
$allTrades.Cells.Item(424, 6).Value = 0.21
$allTrades.Cells.Item(424, 8).Value = "OPEN"
$allTrades.Cells.Item(424, 9).Value = 0
$allTrades.Cells.Item(424, 10).Value = 0
$allTrades.Cells.Item(424, 11).Value = 100.1931846556633
$allTrades.Cells.Item(424, 13).Value = 0
$allTrades.Cells.Item(424, 14).Value = 0
$allTrades.Cells.Item(424, 15).Value = 0
$allTrades.Cells.Item(424, 16).Value = 0.95
$allTrades.Cells.Item(424, 17).Value = "Mean reversion UP: price 18.06% below mean (z=-2.38)"

$allTrades.Cells.Item(425, 1).Value = 424
Set-DateTextCell $allTrades 425 2 "2026-02-18"
$allTrades.Cells.Item(425, 3).Value = "02:03:44"
$allTrades.Cells.Item(425, 4).Value = "MarketMaking"
$allTrades.Cells.Item(425, 5).Value = "DOWN"
$allTrades.Cells.Item(425, 6).Value = 0.79
$allTrades.Cells.Item(425, 8).Value = "OPEN"
$allTrades.Cells.Item(425, 9).Value = 0
$allTrades.Cells.Item(425, 10).Value = 0
$allTrades.Cells.Item(425, 11).Value = 99.43092696038379
$allTrades.Cells.Item(425, 13).Value = 0
$allTrades.Cells.Item(425, 14).Value = 0
$allTrades.Cells.Item(425, 15).Value = 0
$allTrades.Cells.Item(425, 16).Value = 0.6
$allTrades.Cells.Item(425, 17).Value = "Normal spread capture: 247 bps"

$allTrades.Cells.Item(426, 1).Value = 425
Set-DateTextCell $allTrades 426 2 "2026-02-18"
$allTrades.Cells.Item(426, 3).Value = "02:03:46"
$allTrades.Cells.Item(426, 4).Value = "EMAArbitrage"
$allTrades.Cells.Item(426, 5).Value = "DOWN"
$allTrades.Cells.Item(426, 6).Value = 0.72
$allTrades.Cells.Item(426, 8).Value = "OPEN"
$allTrades.Cells.Item(426, 9).Value = 0
$allTrades.Cells.Item(426, 10).Value = 0
$allTrades.Cells.Item(426, 11).Value = 100.450616878256
$allTrades.Cells.Item(426, 13).Value = 0
$allTrades.Cells.Item(426, 14).Value = 0
$allTrades.Cells.Item(426, 15).Value = 0
$allTrades.Cells.Item(426, 16).Value = 0.7059
$allTrades.Cells.Item(426, 17).Value = "EMA:down, RSI:0.0, ROC:-20.59% | 2/3 DOWN"

# ---------------------------------------------------------------------------
# momentum sheet - close trade #393 (row 72) + append trade #422 (row 83)
# Column layout differs from "All Trades":
#   A Trade#  B Date  C Time  D Strategy  E Side  F Entry  G Exit  H Status
#   I P&L%  J P&L$  K Capital After  L Entry Slippage  M Exit Slippage
#   N Confidence  O Entry Reason  P Exit Reason  Q Duration (min)
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Cells.Item(72, 7).Value = 0.74
$momentum.Cells.Item(72, 8).Value = "CLOSED"
$momentum.Cells.Item(72, 9).Value = -7.5
$momentum.Cells.Item(72, 10).Value = -0.06
$momentum.Cells.Item(72, 11).Value = 99.31999999999999
$momentum.Cells.Item(72, 16).Value = "early_exit"
$momentum.Cells.Item(72, 17).Value = 0.36

$momentum.Cells.Item(83, 1).Value = 422
Set-DateTextCell $momentum 83 2 "2026-02-18"
$momentum.Cells.Item(83, 3).Value = "02:03:39"
$momentum.Cells.Item(83, 4).Value = "momentum"
$momentum.Cells.Item(83, 5).Value = "DOWN"
$momentum.Cells.Item(83, 6).Value = 0.8
$momentum.Cells.Item(83, 8).Value = "OPEN"
$momentum.Cells.Item(83, 9).Value = 0
$momentum.Cells.Item(83, 10).Value = 0
$momentum.Cells.Item(83, 11).Value = 99.37699700270591
$momentum.Cells.Item(83, 12).Value = 0
$momentum.Cells.Item(83, 13).Value = 0
$momentum.Cells.Item(83, 14).Value = 0.9
$momentum.Cells.Item(83, 15).Value = "Downward momentum: -20.588% over 10 samples"
$momentum.Cells.Item(83, 17).Value = 0

# ---------------------------------------------------------------------------
# HighProbConvergence sheet - append trade #423 (row 46)
# ---------------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")

$hpc.Cells.Item(46, 1).Value = 423
Set-DateTextCell $hpc 46 2 "2026-02-18"
$hpc.Cells.Item(46, 3).Value = "02:03:40"
$hpc.Cells.Item(46, 4).Value = "HighProbConvergence"
$hpc.Cells.Item(46, 5).Value = "UP"
$hpc.Cells.Item(46, 6).Value = 0.21
$hpc.Cells.Item(46, 8).Value = "OPEN"
$hpc.Cells.Item(46, 9).Value = 0
$hpc.Cells.Item(46, 10).Value = 0
$hpc.Cells.Item(46, 11).Value = 100.1931846556633
$hpc.Cells.Item(46, 12).Value = 0
$hpc.Cells.Item(46, 13).Value = 0
$hpc.Cells.Item(46, 14).Value = 0.95
$hpc.Cells.Item(46, 15).Value = "Mean reversion UP: price 18.06% below mean (z=-2.38)"
$hpc.Cells.Item(46, 17).Value = 0

# ---------------------------------------------------------------------------
# MarketMaking sheet - append trade #424 (row 249)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Cells.Item(249, 1).Value = 424
Set-DateTextCell $mm 249 2 "2026-02-18"
$mm.Cells.Item(249, 3).Value = "02:03:44"
$mm.Cells.Item(249, 4).Value = "MarketMaking"
$mm.Cells.Item(249, 5).Value = "DOWN"
$mm.Cells.Item(249, 6).Value = 0.79
$mm.Cells.Item(249, 8).Value = "OPEN"
$mm.Cells.Item(249, 9).Value = 0
$mm.Cells.Item(249, 10).Value = 0
$mm.Cells.Item(249, 11).Value = 99.43092696038379
$mm.Cells.Item(249, 12).Value = 0
$mm.Cells.Item(249, 13).Value = 0
$mm.Cells.Item(249, 14).Value = 0.6
$mm.Cells.Item(249, 15).Value = "Normal spread capture: 247 bps"
$mm.Cells.Item(249, 17).Value = 0

# ---------------------------------------------------------------------------
# EMAArbitrage sheet - append trade #425 (row 19)
# ---------------------------------------------------------------------------
$ema = $wb.Worksheets.Item("EMAArbitrage")

$ema.Cells.Item(19, 1).Value = 425
Set-DateTextCell $ema 19 2 "2026-02-18"
$ema.Cells.Item(19, 3).Value = "02:03:46"
$ema.Cells.Item(19, 4).Value = "EMAArbitrage"
$ema.Cells.Item(19, 5).Value = "DOWN"
$ema.Cells.Item(19, 6).Value = 0.72
$ema.Cells.Item(19, 8).Value = "OPEN"
$ema.Cells.Item(19, 9).Value = 0
$ema.Cells.Item(19, 10).Value = 0
$ema.Cells.Item(19, 11).Value = 100.450616878256
$ema.Cells.Item(19, 12).Value = 0
$ema.Cells.Item(19, 13).Value = 0
$ema.Cells.Item(19, 14).Value = 0.7059
$ema.Cells.Item(19, 15).Value = "EMA:down, RSI:0.0, ROC:-20.59% | 2/3 DOWN"
$ema.Cells.Item(19, 17).Value = 0
